# Add the new "07/02/2022 - 13/02/2022" reporting week to all three sheets
# of the Toscana school-staff contagion workbook.

$wb = $excel.ActiveWorkbook

$newWeek = "07/02/2022 - 13/02/2022"

# ---------------------------------------------------------------------
# Sheet "Asl Sorveglianza": new rows 117-120
# ---------------------------------------------------------------------
$wsAsl = $wb.Worksheets.Item("Asl Sorveglianza")

$wsAsl.Range("A117").Value = $newWeek
$wsAsl.Range("B117").Value = "AZIENDA USL TOSCANA SUD-EST"
$wsAsl.Range("C117").Value = 21

$wsAsl.Range("A118").Value = $newWeek
$wsAsl.Range("B118").Value = "AZIENDA USL TOSCANA CENTRO"
$wsAsl.Range("C118").Value = 33

$wsAsl.Range("A119").Value = $newWeek
$wsAsl.Range("B119").Value = "AZIENDA USL TOSCANA NORD-OVEST"
$wsAsl.Range("B119").Font.Color = 0
$wsAsl.Range("C119").Value = 25

$wsAsl.Range("B120").Value = "Totale"
$wsAsl.Range("C120").Value = 79

[void]$wsAsl.Activate()
[void]$wsAsl.Range("A117").Select()

# ---------------------------------------------------------------------
# Sheet "Professione": new rows 86-88
# ---------------------------------------------------------------------
$wsProf = $wb.Worksheets.Item("Professione")

$wsProf.Range("A86").Value = $newWeek
$wsProf.Range("B86").Value = "Insegnante"
$wsProf.Range("C86").Value = 72
$wsProf.Range("D86").Value = 2
$wsProf.Range("D86").NumberFormat = "#,##0"

$wsProf.Range("A87").Value = $newWeek
$wsProf.Range("B87").Value = "Personale non docente"
$wsProf.Range("C87").Value = 7

$wsProf.Range("B88").Value = "Totale"
$wsProf.Range("C88").Value = 79
$wsProf.Range("D88").Value = 2
$wsProf.Range("D88").NumberFormat = "#,##0"

[void]$wsProf.Activate()
[void]$wsProf.Range("A86").Select()

# ---------------------------------------------------------------------
# Sheet "Sesso ed età": new rows 159-166
# ---------------------------------------------------------------------
$wsSesso = $wb.Worksheets.Item("Sesso ed età")

$wsSesso.Range("A159").Value = $newWeek
$wsSesso.Range("B159").Value = "19-34"
$wsSesso.Range("C159").Value = "F"
$wsSesso.Range("D159").Value = 13

$wsSesso.Range("A160").Value = $newWeek
$wsSesso.Range("B160").Value = "19-34"
$wsSesso.Range("B160").Font.Color = 0
$wsSesso.Range("C160").Value = "M"
$wsSesso.Range("D160").Value = 2

$wsSesso.Range("A161").Value = $newWeek
$wsSesso.Range("B161").Value = "35-49"
$wsSesso.Range("C161").Value = "F"
$wsSesso.Range("D161").Value = 26

$wsSesso.Range("A162").Value = $newWeek
$wsSesso.Range("B162").Value = "35-49"
$wsSesso.Range("C162").Value = "M"
$wsSesso.Range("D162").Value = 3

$wsSesso.Range("A163").Value = $newWeek
$wsSesso.Range("B163").Value = "50-64"
$wsSesso.Range("C163").Value = "F"
$wsSesso.Range("D163").Value = 28

$wsSesso.Range("A164").Value = $newWeek
$wsSesso.Range("B164").Value = "50-64"
$wsSesso.Range("C164").Value = "M"
$wsSesso.Range("D164").Value = 2

$wsSesso.Range("A165").Value = $newWeek
$wsSesso.Range("B165").Value = "65-79"
$wsSesso.Range("C165").Value = "F"
$wsSesso.Range("D165").Value = 4

$wsSesso.Range("A166").Value = $newWeek
$wsSesso.Range("B166").Value = "65-79"
$wsSesso.Range("C166").Value = "M"
$wsSesso.Range("D166").Value = 1

[void]$wsSesso.Activate()
[void]$wsSesso.Range("D167").Select()
